# Refresh market-derived profit figures across all Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values mirror the latest scheduled-runner price pull; columns H:N are plain data (no formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 21130.8
$ws.Range("I19").Value = 1900
$ws.Range("J19").Value = 25938.5
$ws.Range("K19").Value = 1900
$ws.Range("L19").Value = 25938.5
$ws.Range("M19").Value = -1725
$ws.Range("N19").Value = -26288.5
$ws.Range("H135").Value = 1065.875
$ws.Range("I135").Value = 592.5925999999999
$ws.Range("K135").Value = 5333.3334
$ws.Range("M135").Value = -2798.3334
$ws.Range("H137").Value = 1355.898
$ws.Range("I137").Value = 1222.5161
$ws.Range("J137").Value = 1585.6111
$ws.Range("K137").Value = 3667.5483
$ws.Range("L137").Value = 4756.8333
$ws.Range("M137").Value = -1117.5483
$ws.Range("N137").Value = -9856.8333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 768091.2
$ws.Range("I32").Value = 860344.5600000001
$ws.Range("K32").Value = 860344.5600000001
$ws.Range("M32").Value = -860057.5600000001
$ws.Range("H74").Value = 1325.4412
$ws.Range("I74").Value = 1027.2667
$ws.Range("J74").Value = 1560.8422
$ws.Range("K74").Value = 1027.2667
$ws.Range("L74").Value = 1560.8422
$ws.Range("M74").Value = -153.2666999999999
$ws.Range("N74").Value = -3308.8422
$ws.Range("H77").Value = 1325.4412
$ws.Range("I77").Value = 1027.2667
$ws.Range("J77").Value = 1560.8422
$ws.Range("K77").Value = 5136.3335
$ws.Range("L77").Value = 7804.211
$ws.Range("M77").Value = -768.3334999999997
$ws.Range("N77").Value = -16540.211
$ws.Range("H88").Value = 1677.5714
$ws.Range("I88").Value = 1480.5454
$ws.Range("J88").Value = 2400
$ws.Range("K88").Value = 1480.5454
$ws.Range("L88").Value = 2400
$ws.Range("M88").Value = -1074.5454
$ws.Range("N88").Value = -3212
$ws.Range("H91").Value = 1677.5714
$ws.Range("I91").Value = 1480.5454
$ws.Range("J91").Value = 2400
$ws.Range("K91").Value = 1480.5454
$ws.Range("L91").Value = 2400
$ws.Range("M91").Value = -76.54539999999997
$ws.Range("N91").Value = -5208
$ws.Range("H102").Value = 2700
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -6994
$ws.Range("H132").Value = 2373.8157
$ws.Range("I132").Value = 1814.6492
$ws.Range("J132").Value = 4051.3157
$ws.Range("K132").Value = 5443.9476
$ws.Range("L132").Value = 12153.9471
$ws.Range("M132").Value = -2913.9476
$ws.Range("N132").Value = -17213.9471

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2270.4546
$ws.Range("I86").Value = 2376.28
$ws.Range("J86").Value = 1939.75
$ws.Range("K86").Value = 2376.28
$ws.Range("L86").Value = 1939.75
$ws.Range("M86").Value = -1253.28
$ws.Range("N86").Value = -4185.75
$ws.Range("H89").Value = 2270.4546
$ws.Range("I89").Value = 2376.28
$ws.Range("J89").Value = 1939.75
$ws.Range("K89").Value = 11881.4
$ws.Range("L89").Value = 9698.75
$ws.Range("M89").Value = -6265.400000000001
$ws.Range("N89").Value = -20930.75
$ws.Range("H94").Value = 1298
$ws.Range("I94").Value = 1106.1
$ws.Range("K94").Value = 1106.1
$ws.Range("M94").Value = -655.0999999999999
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 19245
$ws.Range("I93").Value = 9055.5
$ws.Range("J93").Value = 39624
$ws.Range("K93").Value = 9055.5
$ws.Range("L93").Value = 39624
$ws.Range("M93").Value = -7183.5
$ws.Range("N93").Value = -43368
$ws.Range("H94").Value = 1394.6428
$ws.Range("I94").Value = 994.5
$ws.Range("J94").Value = 1461.3334
$ws.Range("K94").Value = 994.5
$ws.Range("L94").Value = 1461.3334
$ws.Range("M94").Value = -543.5
$ws.Range("N94").Value = -2363.3334
$ws.Range("H99").Value = 2081.9092
$ws.Range("J99").Value = 2100
$ws.Range("L99").Value = 2100
$ws.Range("N99").Value = -5096
$ws.Range("H126").Value = 2081.9092
$ws.Range("J126").Value = 2100
$ws.Range("L126").Value = 6300
$ws.Range("N126").Value = -11240
$ws.Range("H132").Value = 2565020.8
$ws.Range("I132").Value = 770.69385
$ws.Range("J132").Value = 10418037
$ws.Range("K132").Value = 2312.08155
$ws.Range("L132").Value = 31254111
$ws.Range("M132").Value = 217.9184500000001
$ws.Range("N132").Value = -31259171

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2843
$ws.Range("J93").Value = 3615.3333
$ws.Range("L93").Value = 10845.9999
$ws.Range("N93").Value = -14589.9999
$ws.Range("H107").Value = 313.10526
$ws.Range("I107").Value = 289.8125
$ws.Range("J107").Value = 437.33334
$ws.Range("K107").Value = 869.4375
$ws.Range("L107").Value = 1312.00002
$ws.Range("M107").Value = 1050.5625
$ws.Range("N107").Value = -5152.000019999999
$ws.Range("H131").Value = 2819.754
$ws.Range("J131").Value = 3490.4119
$ws.Range("L131").Value = 10471.2357
$ws.Range("N131").Value = -20551.2357
$ws.Range("H140").Value = 1490.95
$ws.Range("I140").Value = 997.7037
$ws.Range("K140").Value = 2993.1111
$ws.Range("M140").Value = 2186.8889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8516751
$ws.Range("J10").Value = 275126
$ws.Range("L10").Value = 275126
$ws.Range("N10").Value = -275464
$ws.Range("H70").Value = 5464.9
$ws.Range("I70").Value = 5397.636
$ws.Range("J70").Value = 5649.875
$ws.Range("K70").Value = 5397.636
$ws.Range("L70").Value = 5649.875
$ws.Range("M70").Value = -5127.636
$ws.Range("N70").Value = -6189.875
$ws.Range("H73").Value = 5464.9
$ws.Range("I73").Value = 5397.636
$ws.Range("J73").Value = 5649.875
$ws.Range("K73").Value = 5397.636
$ws.Range("L73").Value = 5649.875
$ws.Range("M73").Value = -4461.636
$ws.Range("N73").Value = -7521.875
$ws.Range("H102").Value = 2367.3125
$ws.Range("I102").Value = 2422.25
$ws.Range("J102").Value = 2312.375
$ws.Range("K102").Value = 2422.25
$ws.Range("L102").Value = 2312.375
$ws.Range("M102").Value = -800.25
$ws.Range("N102").Value = -5556.375
$ws.Range("H113").Value = 85377.164
$ws.Range("I113").Value = 127053.25
$ws.Range("K113").Value = 127053.25
$ws.Range("M113").Value = -124883.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 70007
$ws.Range("J11").Value = 70007
$ws.Range("L11").Value = 70007
$ws.Range("N11").Value = -70287
$ws.Range("H16").Value = 10990133
$ws.Range("I16").Value = 1436.5
$ws.Range("J16").Value = 28572048
$ws.Range("K16").Value = 1436.5
$ws.Range("L16").Value = 28572048
$ws.Range("M16").Value = -1266.5
$ws.Range("N16").Value = -28572388
$ws.Range("H93").Value = 9963.416999999999
$ws.Range("I93").Value = 13400.75
$ws.Range("J93").Value = 3088.75
$ws.Range("K93").Value = 13400.75
$ws.Range("L93").Value = 3088.75
$ws.Range("M93").Value = -12152.75
$ws.Range("N93").Value = -5584.75
$ws.Range("H132").Value = 1967.55
$ws.Range("I132").Value = 1793.9854
$ws.Range("J132").Value = 2336.375
$ws.Range("K132").Value = 5381.956200000001
$ws.Range("L132").Value = 7009.125
$ws.Range("M132").Value = -2851.956200000001
$ws.Range("N132").Value = -12069.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3903
$ws.Range("I122").Value = 3357.8333
$ws.Range("J122").Value = 4993.3335
$ws.Range("K122").Value = 10073.4999
$ws.Range("L122").Value = 14980.0005
$ws.Range("M122").Value = -7623.499899999999
$ws.Range("N122").Value = -19880.0005
$ws.Range("H132").Value = 2605066.5
$ws.Range("I132").Value = 778.5833
$ws.Range("K132").Value = 2335.7499
$ws.Range("M132").Value = 194.2501000000002
